$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$row = 70

$ws.Cells.Item($row, 1).Value = "Partou"
$ws.Cells.Item($row, 2).Value = "BSO Partou Fluitenkruid 10"
$ws.Cells.Item($row, 3).Value = "VGO"

# Store the report date as plain text (matching the rest of column D),
# not as an Excel date serial number.
$ws.Cells.Item($row, 4).NumberFormat = "@"
$ws.Cells.Item($row, 4).Value = "2024-03-07"

$ws.Cells.Item($row, 5).Value = 0
$ws.Cells.Item($row, 6).Value = 0
$ws.Cells.Item($row, 7).Value = 0
$ws.Cells.Item($row, 8).Value = 1
$ws.Cells.Item($row, 9).Value = 0
$ws.Cells.Item($row, 10).Value = 0
